$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update summary header values -------------------------------------
$ws.Range("E11").Value2 = 446726      # VALOR MORA
$ws.Range("C13").Value2 = 4           # Cant. Trabajadores
$ws.Range("F13").Value2 = 8           # Cant. Periodos

# --- 2. Insert 3 extra data rows so the table grows from 5 to 8 rows ------
# (old data rows were 16-20; we insert before the last row, 20, three times
#  so the new table spans rows 16-23, and the old "last row" style/border
#  ends up on the new last row, 23)
$ws.Range("A20:A22").EntireRow.Insert()

# Copy formatting (borders/fonts/number formats) from row 19 (an ordinary
# data row) down into the freshly inserted rows 20-22 so they match the
# rest of the table instead of carrying blank/default formatting.
$ws.Range("B19:J19").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)
$ws.Range("B21:J21").PasteSpecial(-4122)
$ws.Range("B22:J22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 3. Row 16 stays the same worker, only Salario Basico changes ---------
$ws.Range("G16").Value2 = 0

# --- 4. Row 17 becomes a brand-new worker ----------------------------------
$ws.Range("C17").Value2 = "1065897778"
$ws.Range("D17").Value2 = "DIANA CAROLINA ARIAS ROMERO"
$ws.Range("E17").Value2 = "2102"
$ws.Range("F17").Value2 = 35112
$ws.Range("G17").Value2 = 908526

# --- 5. Rows 18 & 19 become another brand-new worker -----------------------
$ws.Range("C18").Value2 = "1047489615"
$ws.Range("D18").Value2 = "LUZ DARY MARQUEZ GUERRERO"
$ws.Range("E18").Value2 = "1912"
$ws.Range("F18").Value2 = 33125
$ws.Range("G18").Value2 = 781242

$ws.Range("C19").Value2 = "1047489615"
$ws.Range("D19").Value2 = "LUZ DARY MARQUEZ GUERRERO"
$ws.Range("E19").Value2 = "1911"
$ws.Range("F19").Value2 = 33125
$ws.Range("G19").Value2 = 781242

# --- 6. Rows 20-23 are Maria Jose Padilla Castillo again, periods reordered
# (rows 20-22 are brand new blank rows created by the Insert above, so every
#  column, including "Tipo Doc" in B, must be written explicitly)
$ws.Range("B20").Value2 = "CC"
$ws.Range("C20").Value2 = "1045309938"
$ws.Range("D20").Value2 = "MARIA JOSE PADILLA CASTILLO"
$ws.Range("E20").Value2 = "2202"
$ws.Range("F20").Value2 = 36341
$ws.Range("G20").Value2 = 908526

$ws.Range("B21").Value2 = "CC"
$ws.Range("C21").Value2 = "1045309938"
$ws.Range("D21").Value2 = "MARIA JOSE PADILLA CASTILLO"
$ws.Range("E21").Value2 = "2201"
$ws.Range("F21").Value2 = 36341
$ws.Range("G21").Value2 = 908526

$ws.Range("B22").Value2 = "CC"
$ws.Range("C22").Value2 = "1045309938"
$ws.Range("D22").Value2 = "MARIA JOSE PADILLA CASTILLO"
$ws.Range("E22").Value2 = "2112"
$ws.Range("F22").Value2 = 36341
$ws.Range("G22").Value2 = 908526

$ws.Range("C23").Value2 = "1045309938"
$ws.Range("D23").Value2 = "MARIA JOSE PADILLA CASTILLO"
$ws.Range("E23").Value2 = "2111"
$ws.Range("F23").Value2 = 36341
$ws.Range("G23").Value2 = 908526

Write-Host "Edit complete"
